# Error Calculations and Plots
# Remove two rows (RM 232 and SC 92) and update several missing/observed
# value cells to reflect the revised imputation pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely.
$ws.Rows.Item(26).Delete()

# After the above deletion, the row that was "SC 92" shifted up to row 27.
# Delete that row as well.
$ws.Rows.Item(27).Delete()

# Now update individual cell values to match the revised dataset.
$ws.Range("E2").Value = -7.2
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = 17.97
$ws.Range("E6").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = ""
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").Value = ""
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F22").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = 16.48
$ws.Range("E24").Value = ""
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("C29").Value = 11.2
$ws.Range("C30").Value = 11.4
$ws.Range("C31").Value = ""
$ws.Range("E31").Value = -8.1
$ws.Range("C32").Value = ""
$ws.Range("E33").Value = -10.7
